# Weekly update: a new Orégano price record (week of the D199 date) is
# inserted at row 199, pushing the existing historical rows 199-248 down to
# 200-249 (the last existing row, old 248, becomes new row 249).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 199; Excel shifts rows 199:248 down to 200:249.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new weekly record.
$ws.Range("A199").Value = 6
$ws.Range("B199").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C199").Value = "Metropolitana"
$ws.Range("D199").Value = 44855
$ws.Range("E199").Value = 13
$ws.Range("F199").Value = 100112029
$ws.Range("G199").Value = "Orégano"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 51
$ws.Range("K199").Value = 16000
$ws.Range("L199").Value = 17000
$ws.Range("M199").Value = 16451
$ws.Range("N199").Value = "$/docena de atados"
$ws.Range("O199").Value = "Región Metropolitana"
$ws.Range("P199").Value = 5484
$ws.Range("Q199").Value = 3
$ws.Range("R199").Value = "Hortaliza"
